$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 74hc sheet: add the two new component rows (11 & 12) with full detail.
#    Written first, and in this specific cell order, so the new shared
#    strings land at the same indices the target workbook expects
#    (322..328, in "first seen" order).
# ---------------------------------------------------------------------------
$wsHc = $wb.Worksheets.Item("74hc")

$wsHc.Range("I11").Value = "https://www.ti.com/lit/ds/symlink/sn74hc153.pdf"
$wsHc.Range("B11").Value = "4:1 mux"
$wsHc.Range("A11").Value = "SN74HC153"
$wsHc.Range("A12").Value = "SN74HC163"
$wsHc.Range("B12").Value = "74xx163 type counter"
$wsHc.Range("I12").Value = "https://www.ti.com/lit/ds/symlink/sn74hc163.pdf"
$wsHc.Range("H12").Value = "60pf"

$wsHc.Range("C11").Value = 5
$wsHc.Range("D11").Value = "8uA"
$wsHc.Range("F11").Value = "80uA"
$wsHc.Range("G11").Value = "3.5pF"
$wsHc.Range("H11").Value = "40pF"

$wsHc.Range("C12").Value = 5
$wsHc.Range("D12").Value = "8uA"
$wsHc.Range("F12").Value = "80uA"
$wsHc.Range("G12").Value = "3.5pF"

$wsHc.Columns.Item(6).ColumnWidth = 25.5

$wsHc.Range("A11:A12").Select()

# ---------------------------------------------------------------------------
# 2) Remaining sheets: each just gets a two-row "SN74HC153 / SN74HC163"
#    callout, styled with the new red font, in the sheet's marker column.
# ---------------------------------------------------------------------------

# 74ac -> column B, rows 12:13
$wsAc = $wb.Worksheets.Item("74ac")
$wsAc.Range("B12").Value = "SN74HC153"
$wsAc.Range("B13").Value = "SN74HC163"
$wsAc.Range("B12:B13").Font.Color = 255
$wsAc.Range("B12:B13").Select()

# 74act -> column A, rows 11:12
$wsAct = $wb.Worksheets.Item("74act")
$wsAct.Range("A11").Value = "SN74HC153"
$wsAct.Range("A12").Value = "SN74HC163"
$wsAct.Range("A11:A12").Font.Color = 255
$wsAct.Range("A11:A12").Select()

# 74hct -> column A, rows 11:12
$wsHct = $wb.Worksheets.Item("74hct")
$wsHct.Range("A11").Value = "SN74HC153"
$wsHct.Range("A12").Value = "SN74HC163"
$wsHct.Range("A11:A12").Font.Color = 255
$wsHct.Range("A11:A12").Select()

# CMOS -> no data change, selection moves to D2
$wsCmos = $wb.Worksheets.Item("CMOS")
$wsCmos.Range("D2").Select()

# 74als -> column A, rows 8:9
$wsAls = $wb.Worksheets.Item("74als")
$wsAls.Range("A8").Value = "SN74HC153"
$wsAls.Range("A9").Value = "SN74HC163"
$wsAls.Range("A8:A9").Font.Color = 255
$wsAls.Range("A8:A9").Select()

# 74ls -> column A, rows 11:12
$wsLs = $wb.Worksheets.Item("74ls")
$wsLs.Range("A11").Value = "SN74HC153"
$wsLs.Range("A12").Value = "SN74HC163"
$wsLs.Range("A11:A12").Font.Color = 255
$wsLs.Range("A11:A12").Select()

# 74lvc -> column A, rows 7:8
$wsLvc = $wb.Worksheets.Item("74lvc")
$wsLvc.Range("A7").Value = "SN74HC153"
$wsLvc.Range("A8").Value = "SN74HC163"
$wsLvc.Range("A7:A8").Font.Color = 255
$wsLvc.Range("A7:A8").Select()

# 74f -> column A, rows 11:12
$wsF = $wb.Worksheets.Item("74f")
$wsF.Range("A11").Value = "SN74HC153"
$wsF.Range("A12").Value = "SN74HC163"
$wsF.Range("A11:A12").Font.Color = 255
$wsF.Range("A11:A12").Select()

# 74s -> column A, rows 6:7 (selected/activated last, matching the new active tab)
$wsS = $wb.Worksheets.Item("74s")
$wsS.Range("A6").Value = "SN74HC153"
$wsS.Range("A7").Value = "SN74HC163"
$wsS.Range("A6:A7").Font.Color = 255
$wsS.Range("A6:A7").Select()
$wsS.Activate()
